$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 3-18 (columns A-F).
# Values are prefixed with a leading apostrophe so Excel stores every
# cell as text (matching the source workbook, which keeps numeric-looking
# values like "166" or "2" as text, not numbers).
$data = @(
    @("112994-030-45-00_R0.pdf", "ROLLED BASE ANGLE", "112994-030-45-01", "166", "ASTM A36", "2"),
    @("112994-030-45-00_R0.pdf", "PL 3/16 x 104 1/4 x 243 7/8", "112994-030-45-02", "1314", "C.S.", "2"),
    @("112994-030-45-00_R0.pdf", "ROLLED BASE ANGLE", "112994-030-45-03", "166", "ASTM A36", "2"),
    @("112994-030-45-00_R0.pdf", "HSS RECT 2 1/2 X 1 1/2 X 1/8 X 104 1/4 LG", "112994-030-45-04", "468", "ASTM A36", "18"),
    @("112994-030-45-00_R0.pdf", "SHEET 16 GAUGE x 24 1/8 x 45", "112994-030-45-12", "327", "C.S.", "18"),
    @("112994-030-45-00_R0.pdf", "SHEET 16 GAUGE x 24 1/8 x 30", "112994-030-46-04", "218", "C.S.", "18"),
    @("112994-030-45-00_R0.pdf", "HHMB 1/2- 13 UNC x 1 3/4 c/w 2 NUTS & FW's", "112994-030-46-05", "N/A", "C.S.", "314"),
    @("112994-030-45-00_R0.pdf", "HHMB 5/8- 11 UNC x 2 1/4 c/w 2 NUTS & FW's", "112994-030-45-13", "N/A", "C.S.", "24"),
    @("112994-030-45-00_R0.pdf", "HHMB 1/2- 13 UNC x 2 c/w 2 NUTS & FW's", "112994-030-45-14", "N/A", "C.S.", "4"),
    @("112994-030-45-00_R0.pdf", "CLEVIS PIN WITH HAIRPIN COTTER 1/2 DIAMETER X 1-1/4 LONG 1 USABLE LENGTH", "112994-030-45-15", "N/A", "C.S.", "10"),
    @("112994-030-45-00_R0.pdf", "FLAT WASHER REGULAR 1/2 GRADE 5 MIN.", "N/A", "N/A", "C.S.", "80"),
    @("112994-030-45-00_R0.pdf", "RD BAR 3/8 X 11-1/8 LG", "112994-030-39-11", "N/A", "C.S.", "2"),
    @("112994-030-45-01_R0.pdf", "L3x2x0.25", "N/A", "83", "C.S.", "1"),
    @("112994-030-45-01_R0.pdf", "FB 1 3/4 X 1/4", "N/A", "0", "C.S.", "2"),
    @("file_name", "description", "part_number", "weight", "material", "quantity"),
    @("112994-030-45-02_R0.pdf", "PLATE STEEL", "112994-030-45-02", "657", "PL 3/16 x 104 1/4 x 243 7/8", "1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = 3 + $i
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = "'" + $row[$c - 1]
        $cell.Style = "Normal"
    }
}
